# New .ttl from Google sheet has been generated.
# The refreshed export drops 4 obsolete vocabulary rows (the old
# "storage format" / "storage temperature" / "derivative type" /
# "freeze thaw cycles of sample storage" entries, id-amr:10175-10178)
# that used to sit right after row 199 ("long-term storage conditions").
# Removing those rows shifts every following row up by 4, which is how
# the regenerated sheet now reads (dimension A1:AA236 -> A1:AA232).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("200:203").Delete()
